$d = $word.ActiveDocument

# The "Thank you" bullet for seco originally reads (as one run):
#   " with initial ideas, helping with programming and testing and, of
#     course, for the awesome SubtitleDownloader API (...), without
#     which SubCentral would not exist."
# It needs to become three separate runs (same en-US formatting) with
# slightly reworded text:
#   1) " with initial "
#   2) "ideas, helping with programming, "
#   3) "testing and, of course, for the awesome SubtitleDownloader API
#       (...), without which SubCentral would not exist."

$old = " with initial ideas, helping with programming and testing and, of course, for the awesome SubtitleDownloader API (http://forum.team-mediaportal.com/mediaportal-plugins-47/subtitledownloader-2-2-a-71651/), without which SubCentral would not exist."

$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start

    # Remove the old run's text, then rebuild it as three runs so each
    # piece becomes its own <w:r> (matching the target XML structure).
    $rng.Text = ""

    $part1 = " with initial "
    $part2 = "ideas, helping with programming, "
    $part3 = "testing and, of course, for the awesome SubtitleDownloader API (http://forum.team-mediaportal.com/mediaportal-plugins-47/subtitledownloader-2-2-a-71651/), without which SubCentral would not exist."

    $r1 = $d.Range($start, $start)
    $r1.InsertAfter($part1)
    $r1.LanguageID = "en-US"

    $r2 = $d.Range($r1.End, $r1.End)
    $r2.InsertAfter($part2)
    $r2.LanguageID = "en-US"

    $r3 = $d.Range($r2.End, $r2.End)
    $r3.InsertAfter($part3)
    $r3.LanguageID = "en-US"
}
